{"js": "// Rename the `<<court>>` Docmosis placeholder to `<<courtName>>`, but only in\n// the \"In the County Court at <<court>>\" heading that lives inside the\n// notice's top table \u2014 the same text also appears later in a plain\n// paragraph and that occurrence must be left untouched (see commit msg:\n// \"Updated manual hearing notices to consistently read from courtName\n// instead of court\").\nconst body = context.document.body;\n\n// `court` is wrapped in its own run (split from the surrounding `<<`/`>>`\n// delimiters), so a case-sensitive, whole-word search cleanly isolates it\n// from words like \"Court\"/\"courts\" elsewhere in the document.\nconst matches = body.search(\"court\", { matchCase: true, matchWholeWord: true });\nmatches.load(\"items\");\nawait context.sync();\n\n// Find the occurrence that sits inside a table (the bold heading cell at the\n// top of the document) rather than the later plain-paragraph mention. Load\n// `parentTableOrNullObject.isNullObject` for every candidate so we can pick\n// the right one once the sync resolves.\nconst infos = matches.items.map((m) => {\n  const t = m.parentTableOrNullObject;\n  t.load(\"isNullObject\");\n  return { range: m, parentTable: t };\n});\nawait context.sync();\n\nlet chosen = null;\nfor (const info of infos) {\n  if (!info.parentTable.isNullObject) {\n    chosen = info.range;\n    break;\n  }\n}\n// Fall back to the first match in document order if nothing was in a table\n// (keeps the script resilient even if the document structure changes).\nif (!chosen && matches.items.length > 0) {\n  chosen = matches.items[0];\n}\n\nif (chosen) {\n  chosen.insertText(\"Name\", \"After\");\n  await context.sync();\n}\n", "ps1": "# Rename the `<<court>>` Docmosis placeholder to `<<courtName>>`, but only in\n# the \"In the County Court at <<court>>\" heading that lives inside the\n# notice's top table \u2014 the same text also appears later in a plain\n# paragraph and that occurrence must be left untouched (see commit msg:\n# \"Updated manual hearing notices to consistently read from courtName\n# instead of court\").\n$d = $word.ActiveDocument\n\n$target = $null\n\n# `court` sits in its own run inside the heading table cell, wrapped by\n# separate `<` `<` ... `>` `>` runs. Scan every table cell (in document\n# order) for a case-sensitive, whole-word \"court\" match so we land on the\n# bold heading occurrence rather than the later plain-paragraph mention.\nforeach ($tbl in $d.Tables) {\n    foreach ($cell in $tbl.Range.Cells) {\n        $cellRange = $cell.Range\n        if ($cellRange.Find.Execute(\"court\", $true, $true)) {\n            $target = $cellRange\n            break\n        }\n    }\n    if ($target -ne $null) { break }\n}\n\n# Fall back to a whole-document search if the document structure ever\n# changes and the placeholder is no longer inside a table.\nif ($target -eq $null) {\n    $bodyRange = $d.Content\n    if ($bodyRange.Find.Execute(\"court\", $true, $true)) {\n        $target = $bodyRange\n    }\n}\n\nif ($target -ne $null) {\n    $target.InsertAfter(\"Name\")\n}\n"}
